$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stale/completed "Notes" comments pending further mixing/implementation comments.
$ws.Range("H13").ClearContents()
$ws.Range("H25").ClearContents()
$ws.Range("H26").ClearContents()
$ws.Range("H28").ClearContents()

# Update the active selection to match the new cursor location.
$ws.Range("J9").Select()
